$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: 1.1.1.1 - harmless 64 -> 63, undetected 30 -> 31
$ws.Range("B3").Value = 63
$ws.Range("E3").Value = 31

# Row 4: 185.220.101.1 - malicious 11 -> 12, undetected 29 -> 28
$ws.Range("C4").Value = 12
$ws.Range("E4").Value = 28
